$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1. Update the "1M resistor" row (row 23): price + ordering link changed
# ------------------------------------------------------------------
$ws1.Range("C23").Value2 = 0.76
$ws1.Range("H23").Value2 = "http://uk.farnell.com/welwyn/cr0805f-1m0fi/resistor-chip-1m-0-125w-1/dp/2341052"

# ------------------------------------------------------------------
# 2. Insert a new row (24) for the new pressure-sensor smoothing capacitor
#    This shifts the summary rows below it down by one.
# ------------------------------------------------------------------
$ws1.Rows.Item(24).Insert()

$ws1.Range("A24").Value2 = "100nF ceramic"
$ws1.Range("B24").Value2 = 1
$ws1.Range("C24").Value2 = 0.204
$ws1.Range("D24").Formula = "=B24*C24"
$ws1.Range("E24").Value2 = "Farnell"
$ws1.Range("F24").Value2 = "0805"
$ws1.Range("G24").Value2 = "Pressure sensor power smoothing"
$ws1.Range("H24").Value2 = "http://uk.farnell.com/avx/08051c104k4z2a/cap-mlcc-x7r-100nf-100v-0805/dp/1833851"

# Fix up the subtotal formula so it now covers the new row too (D2:D25)
$ws1.Range("C27").Formula = "=SUM(D2:D25)"

# ------------------------------------------------------------------
# 3. Highlight (red fill) all the components that have now been ordered
# ------------------------------------------------------------------
$orderedRows = @(2,6,9,10,11,12,14,17,18,19,20,21,22,23)
foreach ($r in $orderedRows) {
    $ws1.Range("A" + $r + ":H" + $r).Interior.Color = 255
}

# ------------------------------------------------------------------
# 4. Update sheet view / active selections
#    - Components sheet becomes the active tab, selection at E29, no frozen topLeftCell
#    - PCB names and values sheet no longer the active tab, selection moves to B8
# ------------------------------------------------------------------
$ws2.Range("B8").Select()
$ws1.Select()
$ws1.Range("E29").Select()
